$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wrapper element for MODS:
#  - C2: "><update type="MODS"> ...   ->   "><datastream type="md_descriptive" operation="update"> ...
#  - Y2: </mods:mods></update></object>   ->   </mods:mods></datastream></object>
$ws.Range("C2").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("Y2").Value = '</mods:mods></datastream></object>'
